# Apply "Added NavBoxes for all ships" edits to the Warheads worksheet.
# This updates the "Vs. strike craft" (column N) values for each weapon
# row, and refreshes the saved view/selection state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warheads")

# Map of row number -> new value for column N
$updates = @{
    2  = 0.6
    3  = 0.6
    4  = 0.6
    6  = 0.55000000000000004
    7  = 0.55000000000000004
    8  = 0.7
    10 = 0.55000000000000004
    11 = 0.55000000000000004
    12 = 0.7
    20 = 0.35
    21 = 0.35
    22 = 0.85
    24 = 0.25
    25 = 0.25
    26 = 0.9
    28 = 0.25
    29 = 0.25
    30 = 0.9
    38 = 0.2
    39 = 0.2
    42 = 0.1
    43 = 0.1
    46 = 0.1
    47 = 0.1
}

foreach ($row in $updates.Keys) {
    $ws.Range("N$row").Value = $updates[$row]
}

# Update the view/selection state that Excel records for the sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N22").Select()
